$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the hyperlink on the "http://www.dti.gov.ph..." source-url cell (A55) ---
# (must happen before the row insert below while the address is still known/stable)
$ws.Range("A55").Hyperlinks.Delete()

# --- Insert a new blank row above the old row 54, pushing the "Source:" block down ---
# New layout right after the insert:
#   A53 Source:
#   A54 (new, blank)
#   A55 Department of Trade and Industry - Micro Small and Medium Enterprises (MSMEs), Statistics.
#   A56 http://www.dti.gov.ph/dti/index.php?p=321   (still styled like a hyperlink)
#   A57 (blank - was old A56)
#   A60 DTIPH
#   A61 Philippines, Department of Trade... (old citation text)
$ws.Rows(54).Insert()

# --- Swap the (now mis-ordered) url text / blank cell so the url ends up on row 57 ---
$ws.Range("A56").Value = ""
$ws.Range("A57").Value = "http://www.dti.gov.ph/dti/index.php?p=321"

# --- Restyle the url cell (A57) to match the plain "source" (italic) look instead of
#     the old blue/underlined hyperlink look ---
$ws.Range("A57").Font.Italic = $true
$ws.Range("A57").Font.Underline = $false
$ws.Range("A57").Font.ColorIndex = -4105

# --- Update the final citation text (old A60, now A61) ---
$ws.Range("A61").Value = "Defined under Small and Medium Enterprise Development (SMED) Council Resolution No. 01 Series of 2003 dated 16 January 2003, Philippines, Department of Trade and Industry. Available at http://www.dti.gov.ph/dti/index.php?p=532"
